$wb = $excel.ActiveWorkbook

# Environment names in the new, desired row order (rows 2-16), shared by both sheets.
$envNames = @(
    "path_planning",
    "attitude_test",
    "altitude_test",
    "attitude_discrete",
    "attitude_continuous",
    "attitude_fragment",
    "altitude_discrete",
    "altitude_continuous",
    "altitude_fragment",
    "attitude_discrete_a2c",
    "attitude_continuous_A2C",
    "attitude_fragment_A2C",
    "altitude_discrete_A2C",
    "altitude_continuous_A2C",
    "altitude_fragment_A2C"
)

$ws1 = $wb.Worksheets.Item("learn")
$ws2 = $wb.Worksheets.Item("run")

for ($i = 0; $i -lt $envNames.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $envNames[$i]
    $ws1.Cells.Item($row, 2).Value = 0
    $ws2.Cells.Item($row, 1).Value = $envNames[$i]
    $ws2.Cells.Item($row, 2).Value = 0
}

# attitude_fragment_A2C row on the "learn" sheet keeps a non-zero iter count.
$ws1.Cells.Item(13, 2).Value = 1

$wb.Save()
